# C5-PowerPoint.pptx edit — Sun, Apr 26, 2020 11:04:47 AM
#
# 1) Slide 6's table switches to a different built-in table style
#    (new style GUID {276B456D-B161-4F76-8C1E-EF0E8DC48A61}).
# 2) The presentation's applied theme changes from "Integral" to the
#    default "Office Theme" colour palette (the deck's Design/theme).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 6 ("SOURCES OF FINANCE" slide)
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{276B456D-B161-4F76-8C1E-EF0E8DC48A61}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colours: Integral -> Office Theme
# ---------------------------------------------------------------------------
function Convert-HexToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (in that order) for the
# built-in "Office Theme" colour scheme.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToOle $officeThemeColors[$i - 1]
}
